$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1277066666666667
$ws.Range("H2").Value = 0.38312
$ws.Range("I2").Value = 0.1657851361976475
$ws.Range("J2").Value = 0.1974405699349423
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08108666666666665
$ws.Range("N2").Value = 0.24326
$ws.Range("O2").Value = 0.02056149724823249
$ws.Range("P2").Value = 0.02939174694341418
$ws.Range("Q2").Value = 0.01035530791111111
$ws.Range("R2").Value = 0.0931977712
$ws.Range("S2").Value = 0.003408790621725776
$ws.Range("T2").Value = 0.005803123267891295
$ws.Range("G3").Value = 0.1277066666666667
$ws.Range("H3").Value = 0.38312
$ws.Range("I3").Value = 0.1657851361976475
$ws.Range("J3").Value = 0.1974405699349423
$ws.Range("O3").Value = 0.02295278771378271
$ws.Range("P3").Value = 0.03280999043916413
$ws.Range("Q3").Value = 0.01155962434666667
$ws.Range("R3").Value = 0.10403661912
$ws.Range("S3").Value = 0.003805231037245155
$ws.Range("T3").Value = 0.006478023211868573
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.1277066666666667
$ws.Range("H4").Value = 0.38312
$ws.Range("I4").Value = 0.1657851361976475
$ws.Range("J4").Value = 0.1974405699349423
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2176356666666667
$ws.Range("N4").Value = 0.652907
$ws.Range("O4").Value = 0.05518681856388939
$ws.Range("P4").Value = 0.07888710565478799
$ws.Range("Q4").Value = 0.02779352553777778
$ws.Range("R4").Value = 0.25014172984
$ws.Range("S4").Value = 0.00914915423192926
$ws.Range("T4").Value = 0.01557551510099935
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.1277066666666667
$ws.Range("H5").Value = 0.38312
$ws.Range("I5").Value = 0.1657851361976475
$ws.Range("J5").Value = 0.1974405699349423
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.5543775
$ws.Range("N5").Value = 7.108755
$ws.Range("O5").Value = 0.9012988964740954
$ws.Range("P5").Value = 0.8589111569626339
$ws.Range("Q5").Value = 0.4539177026
$ws.Range("R5").Value = 2.7235062156
$ws.Range("S5").Value = 0.1494219603067473
$ws.Range("T5").Value = 0.1695839083541831
$ws.Range("A6").Value = "MuSCs"
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.370511
$ws.Range("H6").Value = 0.741022
$ws.Range("I6").Value = 0.4809867660085082
$ws.Range("J6").Value = 0.3818850647690823
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.08108666666666665
$ws.Range("N6").Value = 0.24326
$ws.Range("O6").Value = 0.02056149724823249
$ws.Range("P6").Value = 0.02939174694341418
$ws.Range("Q6").Value = 0.03004350195333333
$ws.Range("R6").Value = 0.18026101172
$ws.Range("S6").Value = 0.009889808065720184
$ws.Range("T6").Value = 0.0112242691851622
$ws.Range("A7").Value = "MuSCs"
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.370511
$ws.Range("H7").Value = 0.741022
$ws.Range("I7").Value = 0.4809867660085082
$ws.Range("J7").Value = 0.3818850647690823
$ws.Range("O7").Value = 0.02295278771378271
$ws.Range("P7").Value = 0.03280999043916413
$ws.Range("Q7").Value = 0.033537544187
$ws.Range("R7").Value = 0.201225265122
$ws.Range("S7").Value = 0.01103998713333216
$ws.Range("T7").Value = 0.01252964532393317
$ws.Range("A8").Value = "MuSCs"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.370511
$ws.Range("H8").Value = 0.741022
$ws.Range("I8").Value = 0.4809867660085082
$ws.Range("J8").Value = 0.3818850647690823
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2176356666666667
$ws.Range("N8").Value = 0.652907
$ws.Range("O8").Value = 0.05518681856388939
$ws.Range("P8").Value = 0.07888710565478799
$ws.Range("Q8").Value = 0.08063640849233333
$ws.Range("R8").Value = 0.483818450954
$ws.Range("S8").Value = 0.02654412938734346
$ws.Range("T8").Value = 0.03012580745242415
$ws.Range("A9").Value = "MuSCs"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.370511
$ws.Range("H9").Value = 0.741022
$ws.Range("I9").Value = 0.4809867660085082
$ws.Range("J9").Value = 0.3818850647690823
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.5543775
$ws.Range("N9").Value = 7.108755
$ws.Range("O9").Value = 0.9012988964740954
$ws.Range("P9").Value = 0.8589111569626339
$ws.Range("Q9").Value = 1.3169359619025
$ws.Range("R9").Value = 5.26774384761
$ws.Range("S9").Value = 0.4335128414221124
$ws.Range("T9").Value = 0.3280053428075629
$ws.Range("G10").Value = 0.2720966666666667
$ws.Range("H10").Value = 0.81629
$ws.Range("I10").Value = 0.3532280977938443
$ws.Range("J10").Value = 0.4206743652959753
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.08108666666666665
$ws.Range("N10").Value = 0.24326
$ws.Range("O10").Value = 0.02056149724823249
$ws.Range("P10").Value = 0.02939174694341418
$ws.Range("Q10").Value = 0.02206341171111111
$ws.Range("R10").Value = 0.1985707054
$ws.Range("S10").Value = 0.007262898560786527
$ws.Range("T10").Value = 0.01236435449036068
$ws.Range("G11").Value = 0.2720966666666667
$ws.Range("H11").Value = 0.81629
$ws.Range("I11").Value = 0.3532280977938443
$ws.Range("J11").Value = 0.4206743652959753
$ws.Range("O11").Value = 0.02295278771378271
$ws.Range("P11").Value = 0.03280999043916413
$ws.Range("Q11").Value = 0.02462937397666667
$ws.Range("R11").Value = 0.22166436579
$ws.Range("S11").Value = 0.008107569543205387
$ws.Range("T11").Value = 0.01380232190336239
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 0.2720966666666667
$ws.Range("H12").Value = 0.81629
$ws.Range("I12").Value = 0.3532280977938443
$ws.Range("J12").Value = 0.4206743652959753
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2176356666666667
$ws.Range("N12").Value = 0.652907
$ws.Range("O12").Value = 0.05518681856388939
$ws.Range("P12").Value = 0.07888710565478799
$ws.Range("Q12").Value = 0.05921793944777778
$ws.Range("R12").Value = 0.53296145503
$ws.Range("S12").Value = 0.01949353494461666
$ws.Range("T12").Value = 0.03318578310136448
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 0.2720966666666667
$ws.Range("H13").Value = 0.81629
$ws.Range("I13").Value = 0.3532280977938443
$ws.Range("J13").Value = 0.4206743652959753
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.5543775
$ws.Range("N13").Value = 7.108755
$ws.Range("O13").Value = 0.9012988964740954
$ws.Range("P13").Value = 0.8589111569626339
$ws.Range("Q13").Value = 0.9671342698250001
$ws.Range("R13").Value = 5.80280561895
$ws.Range("S13").Value = 0.3183640947452357
$ws.Range("T13").Value = 0.3613219058008879
